# Adds a new weekly price record for "Hortaliza, Vega Modelo de Temuco - Alcachofa".
# A brand-new row is inserted at row 100 (pushing the previous rows 100-183 down to
# 101-184), and the new row is populated with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row above the current row 100 - this shifts every existing
# row (100..183) down by one, turning them into rows 101..184 and growing the
# sheet's used range from A1:R183 to A1:R184.
$ws.Rows.Item(100).Insert()

# Populate the newly inserted row 100 with the new weekly record.
$ws.Range("A100").Value = 10
$ws.Range("B100").Value = "Vega Modelo de Temuco"
$ws.Range("C100").Value = "La Araucanía"
$ws.Range("D100").Value = "2022-07-20"
$ws.Range("E100").Value = 9
$ws.Range("F100").Value = 100112013
$ws.Range("G100").Value = "Alcachofa"
$ws.Range("H100").Value = "Española"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 110
$ws.Range("K100").Value = 18000
$ws.Range("L100").Value = 20000
$ws.Range("M100").Value = 18818
$ws.Range("N100").Value = "$/caja 30 unidades"
$ws.Range("O100").Value = "Provincia de Limarí"
$ws.Range("P100").Value = 627
$ws.Range("Q100").Value = 30
$ws.Range("R100").Value = "Hortaliza"
